$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'317.25"
$ws.Range("E2").Value = "'1.79%"
$ws.Range("D3").Value = "'38.06"
$ws.Range("E3").Value = "'0.91%"
$ws.Range("D4").Value = "'5.190"
$ws.Range("E4").Value = "'1.20%"
$ws.Range("D5").Value = "'0.08028"
$ws.Range("E5").Value = "'1.95%"
$ws.Range("D6").Value = "'4.496"
$ws.Range("E6").Value = "'1.92%"
$ws.Range("D7").Value = "'8.511"
$ws.Range("E7").Value = "'2.90%"
$ws.Range("D8").Value = "'1.947"
$ws.Range("E8").Value = "'1.84%"
$ws.Range("E9").Value = "'4.68%"
$ws.Range("D10").Value = "'0.9396"
$ws.Range("E10").Value = "'2.30%"
$ws.Range("D11").Value = "'0.1291"
$ws.Range("E11").Value = "'9.87%"
$ws.Range("E12").Value = "'0.64%"
$ws.Range("D13").Value = "'0.09059"
$ws.Range("E13").Value = "'0.34%"
$ws.Range("D14").Value = "'0.03370"
$ws.Range("E14").Value = "'1.59%"
$ws.Range("D15").Value = "'0.09543"
$ws.Range("E15").Value = "'-0.54%"
$ws.Range("D16").Value = "'0.001406"
$ws.Range("E16").Value = "'1.59%"
$ws.Range("D17").Value = "'0.006397"
$ws.Range("E17").Value = "'9.12%"
$ws.Range("D18").Value = "'3.380"
$ws.Range("E18").Value = "'-4.34%"
$ws.Range("D19").Value = "'0.3523"
$ws.Range("E19").Value = "'2.35%"
$ws.Range("D20").Value = "'6.544"
$ws.Range("E20").Value = "'24.34%"
$ws.Range("E21").Value = "'2.40%"
$ws.Range("D22").Value = "'0.2425"
$ws.Range("E22").Value = "'-6.32%"
$ws.Range("D23").Value = "'0.04384"
$ws.Range("E23").Value = "'0.64%"
$ws.Range("D24").Value = "'0.001230"
$ws.Range("E24").Value = "'-1.44%"
$ws.Range("D25").Value = "'0.004274"
$ws.Range("E25").Value = "'-8.43%"
$ws.Range("D26").Value = "'0.0001330"
$ws.Range("E26").Value = "'-2.15%"
$ws.Range("D27").Value = "'0.0003986"
$ws.Range("E27").Value = "'-0.05%"
$ws.Range("D39").Value = "'0.02363"
$ws.Range("E39").Value = "'5.03%"
$ws.Range("D40").Value = "'0.05164"
$ws.Range("E40").Value = "'1.55%"
$ws.Range("D41").Value = "'0.007654"
$ws.Range("E41").Value = "'2.59%"
$ws.Range("D42").Value = "'0.1404"
$ws.Range("E42").Value = "'3.85%"
$ws.Range("D43").Value = "'0.008635"
$ws.Range("E43").Value = "'-4.40%"
$ws.Range("D44").Value = "'0.002110"
$ws.Range("E44").Value = "'5.57%"
$ws.Range("D45").Value = "'0.008866"
$ws.Range("E45").Value = "'2.77%"
$ws.Range("D46").Value = "'0.00006476"
$ws.Range("E46").Value = "'-1.44%"
$ws.Range("E47").Value = "'-0.06%"
$ws.Range("D48").Value = "'0.002864"
$ws.Range("E48").Value = "'-4.69%"
$ws.Range("D49").Value = "'0.001688"
$ws.Range("E49").Value = "'68.85%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.06%"
